# Auto-generated edit script applying the Garuda_Profits.xlsx market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2031.6364
$ws.Cells.Item(116, 9).Value = 1999.6923
$ws.Cells.Item(116, 10).Value = 2077.7778
$ws.Cells.Item(116, 11).Value = 1999.6923
$ws.Cells.Item(116, 12).Value = 2077.7778
$ws.Cells.Item(116, 13).Value = 1442.3077
$ws.Cells.Item(116, 14).Value = -8961.7778
$ws.Cells.Item(125, 8).Value = 1516.5
$ws.Cells.Item(125, 9).Value = 616
$ws.Cells.Item(125, 10).Value = 1816.6666
$ws.Cells.Item(125, 11).Value = 5544
$ws.Cells.Item(125, 12).Value = 16349.9994
$ws.Cells.Item(125, 13).Value = -3084
$ws.Cells.Item(125, 14).Value = -21269.9994
$ws.Cells.Item(137, 8).Value = 1749.75
$ws.Cells.Item(137, 9).Value = 1617.3572
$ws.Cells.Item(137, 10).Value = 2676.5
$ws.Cells.Item(137, 11).Value = 4852.071599999999
$ws.Cells.Item(137, 12).Value = 8029.5
$ws.Cells.Item(137, 13).Value = -2302.071599999999
$ws.Cells.Item(137, 14).Value = -13129.5
$ws.Cells.Item(141, 8).Value = 2360.2
$ws.Cells.Item(141, 9).Value = 1481.3889
$ws.Cells.Item(141, 10).Value = 2946.074
$ws.Cells.Item(141, 11).Value = 4444.1667
$ws.Cells.Item(141, 12).Value = 8838.222
$ws.Cells.Item(141, 13).Value = 735.8333000000002
$ws.Cells.Item(141, 14).Value = -19198.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3307.1428
$ws.Cells.Item(61, 9).Value = 3380
$ws.Cells.Item(61, 10).Value = 3266.6667
$ws.Cells.Item(61, 11).Value = 3380
$ws.Cells.Item(61, 12).Value = 3266.6667
$ws.Cells.Item(61, 13).Value = -3168
$ws.Cells.Item(61, 14).Value = -3690.6667
$ws.Cells.Item(74, 8).Value = 904.4815
$ws.Cells.Item(74, 9).Value = 964.1667
$ws.Cells.Item(74, 10).Value = 785.1111
$ws.Cells.Item(74, 11).Value = 964.1667
$ws.Cells.Item(74, 12).Value = 785.1111
$ws.Cells.Item(74, 13).Value = -90.16669999999999
$ws.Cells.Item(74, 14).Value = -2533.1111
$ws.Cells.Item(77, 8).Value = 904.4815
$ws.Cells.Item(77, 9).Value = 964.1667
$ws.Cells.Item(77, 10).Value = 785.1111
$ws.Cells.Item(77, 11).Value = 4820.8335
$ws.Cells.Item(77, 12).Value = 3925.5555
$ws.Cells.Item(77, 13).Value = -452.8334999999997
$ws.Cells.Item(77, 14).Value = -12661.5555
$ws.Cells.Item(102, 8).Value = 1886.25
$ws.Cells.Item(102, 9).Value = 1886.25
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1886.25
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -264.25
$ws.Cells.Item(102, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4002.1282
$ws.Cells.Item(132, 9).Value = 5919.6816
$ws.Cells.Item(132, 10).Value = 1520.5883
$ws.Cells.Item(132, 11).Value = 17759.0448
$ws.Cells.Item(132, 12).Value = 4561.7649
$ws.Cells.Item(132, 13).Value = -15229.0448
$ws.Cells.Item(132, 14).Value = -9621.7649
$ws.Cells.Item(136, 8).Value = 3307.1428
$ws.Cells.Item(136, 9).Value = 3380
$ws.Cells.Item(136, 10).Value = 3266.6667
$ws.Cells.Item(136, 11).Value = 10140
$ws.Cells.Item(136, 12).Value = 9800.000100000001
$ws.Cells.Item(136, 13).Value = -7590
$ws.Cells.Item(136, 14).Value = -14900.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 19086.482
$ws.Cells.Item(134, 9).Value = 25407.523
$ws.Cells.Item(134, 10).Value = 2493.75
$ws.Cells.Item(134, 11).Value = 76222.569
$ws.Cells.Item(134, 12).Value = 7481.25
$ws.Cells.Item(134, 13).Value = -73687.569
$ws.Cells.Item(134, 14).Value = -12551.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5716102
$ws.Cells.Item(31, 9).Value = 1418.0714
$ws.Cells.Item(31, 10).Value = 28574838
$ws.Cells.Item(31, 11).Value = 1418.0714
$ws.Cells.Item(31, 12).Value = 28574838
$ws.Cells.Item(31, 13).Value = -1123.0714
$ws.Cells.Item(31, 14).Value = -28575428
$ws.Cells.Item(34, 8).Value = 5716102
$ws.Cells.Item(34, 9).Value = 1418.0714
$ws.Cells.Item(34, 10).Value = 28574838
$ws.Cells.Item(34, 11).Value = 1418.0714
$ws.Cells.Item(34, 12).Value = 28574838
$ws.Cells.Item(34, 13).Value = -1216.0714
$ws.Cells.Item(34, 14).Value = -28575242
$ws.Cells.Item(58, 8).Value = 1015.34375
$ws.Cells.Item(58, 9).Value = 955.2222
$ws.Cells.Item(58, 10).Value = 1340
$ws.Cells.Item(58, 11).Value = 955.2222
$ws.Cells.Item(58, 12).Value = 1340
$ws.Cells.Item(58, 13).Value = -752.2222
$ws.Cells.Item(58, 14).Value = -1746
$ws.Cells.Item(122, 8).Value = 680.93335
$ws.Cells.Item(122, 9).Value = 618.1818
$ws.Cells.Item(122, 10).Value = 853.5
$ws.Cells.Item(122, 11).Value = 1854.5454
$ws.Cells.Item(122, 12).Value = 2560.5
$ws.Cells.Item(122, 13).Value = 595.4546
$ws.Cells.Item(122, 14).Value = -7460.5
$ws.Cells.Item(132, 8).Value = 2088.842
$ws.Cells.Item(132, 9).Value = 2091.0908
$ws.Cells.Item(132, 10).Value = 2085.75
$ws.Cells.Item(132, 11).Value = 6273.2724
$ws.Cells.Item(132, 12).Value = 6257.25
$ws.Cells.Item(132, 13).Value = -3743.2724
$ws.Cells.Item(132, 14).Value = -11317.25
$ws.Cells.Item(134, 8).Value = 1161.625
$ws.Cells.Item(134, 9).Value = 1061.3334
$ws.Cells.Item(134, 11).Value = 3184.0002
$ws.Cells.Item(134, 13).Value = -649.0001999999999
$ws.Cells.Item(136, 8).Value = 1015.34375
$ws.Cells.Item(136, 9).Value = 955.2222
$ws.Cells.Item(136, 10).Value = 1340
$ws.Cells.Item(136, 11).Value = 2865.6666
$ws.Cells.Item(136, 12).Value = 4020
$ws.Cells.Item(136, 13).Value = -315.6666
$ws.Cells.Item(136, 14).Value = -9120

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 823.7
$ws.Cells.Item(5, 9).Value = 408
$ws.Cells.Item(5, 10).Value = 1239.4
$ws.Cells.Item(5, 11).Value = 1224
$ws.Cells.Item(5, 12).Value = 3718.2
$ws.Cells.Item(5, 13).Value = -1112
$ws.Cells.Item(5, 14).Value = -3942.2
$ws.Cells.Item(131, 8).Value = 4119981.8
$ws.Cells.Item(131, 9).Value = 16111.429
$ws.Cells.Item(131, 10).Value = 5556336.5
$ws.Cells.Item(131, 11).Value = 48334.287
$ws.Cells.Item(131, 12).Value = 16669009.5
$ws.Cells.Item(131, 13).Value = -43294.287
$ws.Cells.Item(131, 14).Value = -16679089.5
$ws.Cells.Item(135, 8).Value = 823.7
$ws.Cells.Item(135, 9).Value = 408
$ws.Cells.Item(135, 10).Value = 1239.4
$ws.Cells.Item(135, 11).Value = 3672
$ws.Cells.Item(135, 12).Value = 11154.6
$ws.Cells.Item(135, 13).Value = -1137
$ws.Cells.Item(135, 14).Value = -16224.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 76266.63
$ws.Cells.Item(132, 9).Value = 119129.47
$ws.Cells.Item(132, 11).Value = 357388.41
$ws.Cells.Item(132, 13).Value = -354858.41

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1427.2727
$ws.Cells.Item(61, 9).Value = 1344.4445
$ws.Cells.Item(61, 10).Value = 1800
$ws.Cells.Item(61, 11).Value = 1344.4445
$ws.Cells.Item(61, 12).Value = 1800
$ws.Cells.Item(61, 13).Value = -1142.4445
$ws.Cells.Item(61, 14).Value = -2204
$ws.Cells.Item(100, 8).Value = 670.58826
$ws.Cells.Item(100, 9).Value = 1000
$ws.Cells.Item(100, 11).Value = 1000
$ws.Cells.Item(100, 13).Value = -459
$ws.Cells.Item(113, 8).Value = 1427.2727
$ws.Cells.Item(113, 9).Value = 1344.4445
$ws.Cells.Item(113, 10).Value = 1800
$ws.Cells.Item(113, 11).Value = 1344.4445
$ws.Cells.Item(113, 12).Value = 1800
$ws.Cells.Item(113, 13).Value = 825.5554999999999
$ws.Cells.Item(113, 14).Value = -6140
$ws.Cells.Item(132, 8).Value = 6418.512
$ws.Cells.Item(132, 9).Value = 7415.4194
$ws.Cells.Item(132, 10).Value = 3328.1
$ws.Cells.Item(132, 11).Value = 22246.2582
$ws.Cells.Item(132, 12).Value = 9984.299999999999
$ws.Cells.Item(132, 13).Value = -19716.2582
$ws.Cells.Item(132, 14).Value = -15044.3
$ws.Cells.Item(136, 8).Value = 10176
$ws.Cells.Item(136, 9).Value = 15301.714
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 45905.142
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -43355.142
$ws.Cells.Item(136, 14).Value = -14100
$ws.Cells.Item(140, 8).Value = 35000
$ws.Cells.Item(140, 10).Value = 35000
$ws.Cells.Item(140, 12).Value = 35000
$ws.Cells.Item(140, 14).Value = -45360

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1982.88
$ws.Cells.Item(132, 9).Value = 1503.7
$ws.Cells.Item(132, 10).Value = 3899.6
$ws.Cells.Item(132, 11).Value = 4511.1
$ws.Cells.Item(132, 12).Value = 11698.8
$ws.Cells.Item(132, 13).Value = -1981.1
$ws.Cells.Item(132, 14).Value = -16758.8
$ws.Cells.Item(136, 8).Value = 6218.1304
$ws.Cells.Item(136, 9).Value = 8610.4375
$ws.Cells.Item(136, 10).Value = 750
$ws.Cells.Item(136, 11).Value = 25831.3125
$ws.Cells.Item(136, 12).Value = 2250
$ws.Cells.Item(136, 13).Value = -23281.3125
$ws.Cells.Item(136, 14).Value = -7350
$ws.Cells.Item(138, 8).Value = 96714.5
$ws.Cells.Item(138, 10).Value = 96714.5
$ws.Cells.Item(138, 12).Value = 96714.5
$ws.Cells.Item(138, 14).Value = -106994.5
